$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 658
$ws.Range("J17").Value = 658
$ws.Range("L17").Value = 1974
$ws.Range("N17").Value = -2310
$ws.Range("H32").Value = 1499.6666
$ws.Range("J32").Value = 1999
$ws.Range("L32").Value = 1999
$ws.Range("N32").Value = -2651
$ws.Range("H87").Value = 123984.5
$ws.Range("J87").Value = 123984.5
$ws.Range("L87").Value = 123984.5
$ws.Range("N87").Value = -126480.5
$ws.Range("H90").Value = 123984.5
$ws.Range("J90").Value = 123984.5
$ws.Range("L90").Value = 371953.5
$ws.Range("N90").Value = -384433.5
$ws.Range("H92").Value = 576.4091
$ws.Range("I92").Value = 565.0526
$ws.Range("J92").Value = 648.3333
$ws.Range("K92").Value = 565.0526
$ws.Range("L92").Value = 648.3333
$ws.Range("M92").Value = 682.9474
$ws.Range("N92").Value = -3144.3333
$ws.Range("H93").Value = 98291.664
$ws.Range("J93").Value = 98291.664
$ws.Range("L93").Value = 98291.664
$ws.Range("N93").Value = -103283.664
$ws.Range("H94").Value = 3005
$ws.Range("I94").Value = 3005
$ws.Range("K94").Value = 3005
$ws.Range("M94").Value = -2554
$ws.Range("H101").Value = 2884.2856
$ws.Range("I101").Value = 3015
$ws.Range("K101").Value = 9045
$ws.Range("M101").Value = -7423
$ws.Range("H113").Value = 71435784
$ws.Range("J113").Value = 66675084
$ws.Range("L113").Value = 66675084
$ws.Range("N113").Value = -66681592
$ws.Range("H132").Value = 2405.0952
$ws.Range("I132").Value = 1986.7222
$ws.Range("J132").Value = 4915.3335
$ws.Range("K132").Value = 5960.1666
$ws.Range("L132").Value = 14746.0005
$ws.Range("M132").Value = -3430.1666
$ws.Range("N132").Value = -19806.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1759.3636
$ws.Range("I2").Value = 1685.4
$ws.Range("J2").Value = 2499
$ws.Range("K2").Value = 1685.4
$ws.Range("L2").Value = 2499
$ws.Range("M2").Value = -1572.4
$ws.Range("N2").Value = -2725
$ws.Range("H45").Value = 2647.1177
$ws.Range("I45").Value = 2370.2
$ws.Range("K45").Value = 2370.2
$ws.Range("M45").Value = -1993.2
$ws.Range("H116").Value = 1759.3636
$ws.Range("I116").Value = 1685.4
$ws.Range("J116").Value = 2499
$ws.Range("K116").Value = 1685.4
$ws.Range("L116").Value = 2499
$ws.Range("M116").Value = 608.5999999999999
$ws.Range("N116").Value = -7087

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1759.3636
$ws.Range("I3").Value = 1685.4
$ws.Range("J3").Value = 2499
$ws.Range("K3").Value = 1685.4
$ws.Range("L3").Value = 2499
$ws.Range("M3").Value = -1571.4
$ws.Range("N3").Value = -2727
$ws.Range("H99").Value = 2343.6924
$ws.Range("I99").Value = 1768.7142
$ws.Range("K99").Value = 1768.7142
$ws.Range("M99").Value = -270.7141999999999
$ws.Range("H134").Value = 51771.906
$ws.Range("I134").Value = 2144.4285
$ws.Range("K134").Value = 6433.2855
$ws.Range("M134").Value = -3898.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1318.2
$ws.Range("I107").Value = 1252.5
$ws.Range("K107").Value = 1252.5
$ws.Range("M107").Value = 667.5
$ws.Range("H132").Value = 3179.6667
$ws.Range("I132").Value = 3106.6
$ws.Range("K132").Value = 9319.799999999999
$ws.Range("M132").Value = -6789.799999999999
$ws.Range("H134").Value = 630712.25
$ws.Range("I134").Value = 1002048.7
$ws.Range("J134").Value = 11818.167
$ws.Range("K134").Value = 3006146.1
$ws.Range("L134").Value = 35454.501
$ws.Range("M134").Value = -3003611.1
$ws.Range("N134").Value = -40524.501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3604205.8
$ws.Range("I4").Value = 3782593.5
$ws.Range("K4").Value = 11347780.5
$ws.Range("M4").Value = -11347668.5
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2500
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 7500
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -7331
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 2500
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 7500
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -7398
$ws.Range("N27").ClearContents()
$ws.Range("H44").Value = 3998.9167
$ws.Range("I44").Value = 2943
$ws.Range("K44").Value = 8829
$ws.Range("M44").Value = -8431
$ws.Range("H114").Value = 1304
$ws.Range("I114").Value = 907
$ws.Range("J114").Value = 1899.5
$ws.Range("K114").Value = 2721
$ws.Range("L114").Value = 5698.5
$ws.Range("M114").Value = 533
$ws.Range("N114").Value = -12206.5
$ws.Range("H121").Value = 4758
$ws.Range("J121").Value = 10995
$ws.Range("L121").Value = 32985
$ws.Range("N121").Value = -35605
$ws.Range("H122").Value = 1321.1052
$ws.Range("I122").Value = 782.5
$ws.Range("J122").Value = 1569.6923
$ws.Range("K122").Value = 7042.5
$ws.Range("L122").Value = 14127.2307
$ws.Range("M122").Value = -4592.5
$ws.Range("N122").Value = -19027.2307
$ws.Range("H131").Value = 4450.6816
$ws.Range("J131").Value = 2697.3845
$ws.Range("L131").Value = 8092.1535
$ws.Range("N131").Value = -18172.1535
$ws.Range("H134").Value = 4055.5293
$ws.Range("I134").Value = 2412.1667
$ws.Range("K134").Value = 7236.500100000001
$ws.Range("M134").Value = -2166.500100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6303
$ws.Range("I102").Value = 5125.9287
$ws.Range("J102").Value = 7950.9
$ws.Range("K102").Value = 5125.9287
$ws.Range("L102").Value = 7950.9
$ws.Range("M102").Value = -3503.9287
$ws.Range("N102").Value = -11194.9
$ws.Range("H136").Value = 31389.2
$ws.Range("J136").Value = 31389.2
$ws.Range("L136").Value = 94167.60000000001
$ws.Range("N136").Value = -99267.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2076.5881
$ws.Range("I16").Value = 1759.4166
$ws.Range("J16").Value = 2837.8
$ws.Range("K16").Value = 1759.4166
$ws.Range("L16").Value = 2837.8
$ws.Range("M16").Value = -1589.4166
$ws.Range("N16").Value = -3177.8
$ws.Range("H22").Value = 2605.16
$ws.Range("I22").Value = 2590.5
$ws.Range("K22").Value = 2590.5
$ws.Range("M22").Value = -2295.5
$ws.Range("H27").Value = 2605.16
$ws.Range("I27").Value = 2590.5
$ws.Range("K27").Value = 2590.5
$ws.Range("M27").Value = -2483.5
$ws.Range("H55").Value = 58823784
$ws.Range("J55").Value = 224.4
$ws.Range("L55").Value = 224.4
$ws.Range("N55").Value = -570.4
$ws.Range("H82").Value = 718.05884
$ws.Range("J82").Value = 1050.4286
$ws.Range("L82").Value = 1050.4286
$ws.Range("N82").Value = -1772.4286
$ws.Range("H85").Value = 718.05884
$ws.Range("J85").Value = 1050.4286
$ws.Range("L85").Value = 1050.4286
$ws.Range("N85").Value = -3546.4286
$ws.Range("H104").Value = 47185.25
$ws.Range("J104").Value = 47185.25
$ws.Range("L104").Value = 47185.25
$ws.Range("N104").Value = -54173.25
$ws.Range("H122").Value = 5352.846
$ws.Range("I122").Value = 4426.6113
$ws.Range("K122").Value = 13279.8339
$ws.Range("M122").Value = -10829.8339
$ws.Range("H132").Value = 745105.1
$ws.Range("I132").Value = 1112563.5
$ws.Range("K132").Value = 3337690.5
$ws.Range("M132").Value = -3335160.5
$ws.Range("H136").Value = 76951.87
$ws.Range("I136").Value = 3115.7273
$ws.Range("J136").Value = 280001.25
$ws.Range("K136").Value = 9347.1819
$ws.Range("L136").Value = 840003.75
$ws.Range("M136").Value = -6797.1819
$ws.Range("N136").Value = -845103.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5499.5
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -12746
$ws.Range("H107").Value = 26317082
$ws.Range("I107").Value = 35715540
$ws.Range("J107").Value = 1397.8
$ws.Range("K107").Value = 107146620
$ws.Range("L107").Value = 4193.4
$ws.Range("M107").Value = -107144700
$ws.Range("N107").Value = -8033.4
$ws.Range("H126").Value = 912.44446
$ws.Range("I126").Value = 912.44446
$ws.Range("K126").Value = 2737.33338
$ws.Range("M126").Value = -267.33338

Write-Host "Edit complete"